# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Rows 36/37 (ARBITRUM <-> HuobiToken) swap order in addition to value updates.
# Values that look purely numeric ("223.12", "0.5295", ...) are written with a
# leading apostrophe (forces text, matching the source inlineStr cells) and
# then the cell style is reset to "Normal" so no stray NumberFormat/quotePrefix
# styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.238.66"
$ws.Range("E2").Value = "  -2.29%  "

$ws.Range("D3").Value = "1.705.42"
$ws.Range("E3").Value = "  -1.55%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "'223.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.69%  "

$ws.Range("D6").Value = "'0.5295"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "'0.2641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.68%  "

$ws.Range("D9").Value = "'0.06548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.40%  "

$ws.Range("D10").Value = "'20.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.03%  "

$ws.Range("D11").Value = "'0.07649"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "

$ws.Range("D12").Value = "'4.571"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("D13").Value = "1.736.01"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "1.943.31"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "'0.5726"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.87%  "

$ws.Range("D16").Value = "0.0₅8176"
$ws.Range("E16").Value = "  -2.60%  "

$ws.Range("D17").Value = "'67.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").Value = "27.257.68"
$ws.Range("E18").Value = "  -2.11%  "

$ws.Range("D19").Value = "'215.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "'4.662"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.44%  "

$ws.Range("D22").Value = "'10.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.38%  "

$ws.Range("D23").Value = "'5.957"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.69%  "

$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").Value = "'142.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "

$ws.Range("D26").Value = "'1.750"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.63%  "

$ws.Range("E27").Value = "  -2.62%  "

$ws.Range("D28").Value = "'7.253"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.54%  "

$ws.Range("D29").Value = "'16.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.42%  "

$ws.Range("D30").Value = "'0.05365"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.30%  "

$ws.Range("D31").Value = "'1.289"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").Value = "'3.494"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.43%  "

$ws.Range("D33").Value = "'3.413"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("D34").Value = "'1.636"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "'2.877"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9473"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").Value = "'0.5853"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").Value = "'0.01629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("D40").Value = "'5.875"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D41").Value = "'1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").Value = "1.039.07"
$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("D43").Value = "'0.8399"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").Value = "'100.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.90%  "

$ws.Range("D45").Value = "1.849.14"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +6.24%  "

$ws.Range("D47").Value = "'57.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("D48").Value = "'0.4492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("D50").Value = "'8.045"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("D51").Value = "'0.06549"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.13%  "
